$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("uma-title")

# Correct a mistranscribed katakana character in an existing title
$ws.Range("A9").Value = "シューティンスタァ・ルヴゥ"

# Append 5 new uma-title rows (insert rows so they pick up the same
# cell style as the rest of the table, then populate)
for ($i = 48; $i -le 52; $i++) {
    $ws.Rows.Item($i).Insert()
}

$ws.Range("A48").Value = "ボーノ☆アラモーダ"
$ws.Range("B48").Value = "Bono ☆ Alla Moda"

$ws.Range("A49").Value = "超特急！フルカラー特殊PP"
$ws.Range("B49").Value = "Super Express! Full Color Special"

$ws.Range("A50").Value = "Make up Vampire!"
$ws.Range("B50").Value = "Make up Vampire!"

$ws.Range("A51").Value = "シフォンリボンマミー"
$ws.Range("B51").Value = "Chiffon Ribbon Mummy"

$ws.Range("A52").Value = "プリンセス・オブ・ピンク"
$ws.Range("B52").Value = "Princess of Pink"

# Grow the query table to cover the new rows
$lo = $ws.ListObjects.Item("uma_title")
$lo.Resize($ws.Range("A1:B52"))

# Keep the ExternalData_1 defined name (backing range for the query
# table) in sync with the new extent
$name = $wb.Names.Item("ExternalData_1")
$name.RefersTo = "='uma-title'!`$A`$1:`$B`$52"

# Leave the selection where the edits were made
$ws.Activate() | Out-Null
$ws.Range("E30").Select() | Out-Null
